# Generate Report for Handoff
# Refresh the "Ready for handoff" / "Handback transform failed" timestamps
# that the report generator stamps when it regenerates the handback status
# report. The zh-cn and de-de detail sheets get their own freshly-generated
# "Latest Handoff Datetime" (column E) for every row still awaiting/failed
# handoff, and the Overview sheet's "Latest Handoff Date" (column D) is
# refreshed to the newest of the two locale timestamps for those same rows.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Rows whose status is "Handback transform failed" (7) or "Ready for
# handoff" (10-16) pick up the newly generated handoff timestamps.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$zhcnStamp = "2016-03-24 18:36:46"
$dedeStamp = "2016-03-24 18:36:51"
# Overview shows the latest handoff date across locales - de-de is the
# most recent of the two new stamps.
$overviewStamp = "2016-03-24 18:36:51"

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = $zhcnStamp
    $dede.Range("E$r").Value = $dedeStamp
    $overview.Range("D$r").Value = $overviewStamp
}
